# "Las pelotas se generan, se mueven en la pantalla y rebotan.
#  Se detecta la colision con la nave" -- add the two new requirement
# rows to the "Must Have" sheet (D6 / E6), give them the same
# highlight-style formatting used elsewhere in the sheet, and update the
# view/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Must Have")

# New content for D6 / E6 (becomes shared-string entries 62 / 63)
$ws.Range("D6").Value = "Ver si hay colision con alguna pelota"
$ws.Range("E6").Value = "Crear arreglo de pelotas"

# Style the two new cells like the rest of the "checklist" columns:
# centered, filled with the light-green accent used for highlighting.
# -4108 = xlCenter; 11854021 = RGB(197,224,180), the Accent6 "Lighter 60%" swatch
$rng = $ws.Range("D6:E6")
$rng.HorizontalAlignment = -4108
$rng.Interior.Color = 11854021

# Conditional formatting: mirror the "x" highlight rules that already
# exist on column G, applied individually to the new D6/E6 cells
# (green = "x" met, red = alternate state).
$e6 = $ws.Range("E6")
$fcE1 = $e6.FormatConditions.Add(1, 3, '"x"')
$fcE1.Font.Color = 24832
$fcE1.Interior.Color = 13561798

$fcE2 = $e6.FormatConditions.Add(1, 3, '"x"')
$fcE2.Font.Color = 393372
$fcE2.Interior.Color = 13551615

$d6 = $ws.Range("D6")
$fcD1 = $d6.FormatConditions.Add(1, 3, '"x"')
$fcD1.Font.Color = 24832
$fcD1.Interior.Color = 13561798

$fcD2 = $d6.FormatConditions.Add(1, 3, '"x"')
$fcD2.Font.Color = 393372
$fcD2.Interior.Color = 13551615

# Column D grew wider to fit the longer text; column E now matches the
# previous width of column D.
$ws.Columns.Item(4).ColumnWidth = 32.6
$ws.Columns.Item(5).ColumnWidth = 21.6

# Update the view: scroll/selection moved to F7.
$ws.Activate()
$ws.Range("F7").Select()

Write-Output "Must Have sheet updated: D6/E6 populated and formatted"
